# NaukariVisualsData.xlsx edit
# - "Data" sheet: Runmode column (A) bulk-changed from "N" to "Y" for rows 4-13
# - "Data" sheet: new "Total Jobs" figures filled into column F (rows 3-13),
#   with the header F2 given an integer number format, and the data cells
#   given word-wrap.
# - Both sheets still point their "OnlyPositions" cell at the shared string,
#   which is renumbered automatically once the now-unused "N" string is
#   dropped by saving.
# - View state (selections) updated on both sheets.

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("Testcases")
$wsData = $wb.Worksheets.Item("Data")

# ---- Data sheet: Runmode column N -> Y for the remaining position rows ----
$wsData.Range("A4").Value  = "Y"
$wsData.Range("A5").Value  = "Y"
$wsData.Range("A6").Value  = "Y"
$wsData.Range("A7").Value  = "Y"
$wsData.Range("A8").Value  = "Y"
$wsData.Range("A9").Value  = "Y"
$wsData.Range("A10").Value = "Y"
$wsData.Range("A11").Value = "Y"
$wsData.Range("A12").Value = "Y"
$wsData.Range("A13").Value = "Y"

# ---- Data sheet: header "Total Jobs" cell gets an integer display format ----
$wsData.Range("F2").NumberFormat = "0"

# ---- Data sheet: the new "Total Jobs" data cells lose their box border and
#      gain word-wrap ----
$wsData.Range("F3:F13").Borders.LineStyle = -4142
$wsData.Range("F3:F13").WrapText = $true

# ---- Data sheet: fill in the "Total Jobs" counts (written as text, in row
#      order, matching how the source tool produced them) ----
$wsData.Range("F3").Value  = "'4065"
$wsData.Range("F4").Value  = "'638"
$wsData.Range("F5").Value  = "'100192"
$wsData.Range("F6").Value  = "'24498"
$wsData.Range("F7").Value  = "'2887"
$wsData.Range("F8").Value  = "'11697"
$wsData.Range("F9").Value  = "'2531"
$wsData.Range("F10").Value = "'1932"
$wsData.Range("F11").Value = "'309"
$wsData.Range("F12").Value = "'2776"
$wsData.Range("F13").Value = "'7196"

# ---- View state ----
$wsTest.Range("A1:B3").Select() | Out-Null
$wsData.Activate()
$wsData.Range("C21").Select() | Out-Null
